# Fruta / hortaliza, semanal
# Re-shuffle the weekly sample rows (3-46): each row's Fecha/Variedad/
# Volumen/Precio columns are reassigned from another row of the same
# original data set (row 2 and the header stay untouched).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns that carry the per-row sample data affected by the reshuffle.
$cols = @("D", "H", "J", "K", "L", "M", "P")

$firstRow = 3
$lastRow = 46

# Snapshot the current ("before") values for every affected column/row
# so the re-assignment below can be computed from a single consistent
# view of the original data (a simultaneous permutation).
$snapshot = @{}
for ($r = $firstRow; $r -le $lastRow; $r++) {
    $rowData = @{}
    foreach ($c in $cols) {
        $rowData[$c] = $ws.Range("$c$r").Value2
    }
    $snapshot[$r] = $rowData
}

# Mapping: new row number -> source row number (source row's original
# values are copied into the new row).
$rowMap = @{
    3 = 37; 4 = 30; 5 = 7; 6 = 38; 7 = 22; 8 = 8; 9 = 13; 10 = 34;
    11 = 17; 12 = 21; 13 = 26; 14 = 4; 15 = 35; 16 = 11; 17 = 25;
    18 = 41; 19 = 40; 20 = 33; 21 = 15; 22 = 6; 23 = 16; 24 = 12;
    25 = 27; 26 = 10; 27 = 9; 28 = 31; 29 = 46; 30 = 39; 31 = 5;
    32 = 43; 33 = 36; 34 = 23; 35 = 28; 36 = 44; 37 = 29; 38 = 18;
    39 = 14; 40 = 3; 41 = 32; 42 = 24; 43 = 45; 44 = 20; 45 = 42;
    46 = 19
}

foreach ($newRow in ($rowMap.Keys | Sort-Object)) {
    $srcRow = $rowMap[$newRow]
    $srcData = $snapshot[$srcRow]
    foreach ($c in $cols) {
        $ws.Range("$c$newRow").Value = $srcData[$c]
    }
}
